# Insert a new weekly record at row 379, shifting the existing rows
# 379-474 down to 380-475 (same as the author's "Fruta / hortaliza, semanal"
# weekly-refresh commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 379..474 down by one row, preserving their existing content.
$ws.Rows.Item(379).EntireRow.Insert()

# Populate the freshly inserted row 379 with the new weekly entry.
$ws.Range("A379").Value = 9
$ws.Range("B379").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C379").Value = "Metropolitana"
$ws.Range("D379").Value = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(44855)
$ws.Range("E379").Value = 13
$ws.Range("F379").Value = 100112039
$ws.Range("G379").Value = "Ciboulette"
$ws.Range("H379").Value = "Sin especificar"
$ws.Range("I379").Value = "Primera"
$ws.Range("J379").Value = 750
$ws.Range("K379").Value = 800
$ws.Range("L379").Value = 1000
$ws.Range("M379").Value = 920
$ws.Range("N379").Value = "$/docena de atados"
$ws.Range("O379").Value = "Provincia de Chacabuco"
$ws.Range("P379").Value = 307
$ws.Range("Q379").Value = 3
$ws.Range("R379").Value = "Hortaliza"
